# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets,
# and the single changed value on the 演出 sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4701
$ws1.Range("F3").Value = 1861
$ws1.Range("F4").Value = 146
$ws1.Range("F6").Value = 3155
$ws1.Range("F7").Value = 583
$ws1.Range("F9").Value = 280
$ws1.Range("F10").Value = 644
$ws1.Range("F11").Value = 550
$ws1.Range("F13").Value = 401
$ws1.Range("F14").Value = 139
$ws1.Range("F15").Value = 1794
$ws1.Range("F16").Value = 1369
$ws1.Range("F18").Value = 1637
$ws1.Range("F19").Value = 20
$ws1.Range("F20").Value = 129
$ws1.Range("F21").Value = 612
$ws1.Range("F22").Value = 17
$ws1.Range("F30").Value = 37
$ws1.Range("F32").Value = 3960
$ws1.Range("F33").Value = 8
$ws1.Range("F34").Value = 778
$ws1.Range("F36").Value = 1302
$ws1.Range("F38").Value = 1879

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 54

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4701
$ws4.Range("F3").Value = 1861
$ws4.Range("F4").Value = 146
$ws4.Range("F6").Value = 3155
$ws4.Range("F7").Value = 583
$ws4.Range("F9").Value = 280
$ws4.Range("F10").Value = 644
$ws4.Range("F11").Value = 550
$ws4.Range("F14").Value = 401
$ws4.Range("F15").Value = 139
$ws4.Range("F16").Value = 1794
$ws4.Range("F17").Value = 1369
$ws4.Range("F19").Value = 1637
$ws4.Range("F20").Value = 20
$ws4.Range("F21").Value = 129
$ws4.Range("F22").Value = 612
$ws4.Range("F23").Value = 17
$ws4.Range("F31").Value = 37
$ws4.Range("F33").Value = 3960
$ws4.Range("F34").Value = 54
$ws4.Range("F35").Value = 8
$ws4.Range("F37").Value = 778
$ws4.Range("F39").Value = 1302
$ws4.Range("F41").Value = 1879

$wb.Save()
